# Añadimos a StreamLit la pantalla de graficos de analisis de partidos
# Appends the next block of match events (rows 84-109) to the event log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(84, 1).Value = '00:04:43.530'
$ws.Cells.Item(84, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(84, 3).Value = 'Pass'
$ws.Cells.Item(84, 4).Value = 4
$ws.Cells.Item(84, 5).Value = 50
$ws.Cells.Item(84, 6).Value = 60
$ws.Cells.Item(84, 7).Value = 33
$ws.Cells.Item(84, 8).Value = 74
$ws.Cells.Item(84, 9).Value = 'Adrian Pombo'
$ws.Cells.Item(84, 10).Value = 'Complete'
$ws.Cells.Item(84, 11).Value = 'Ground Pass'
$ws.Cells.Item(84, 12).Value = 'Santiago Sanchez'
$ws.Cells.Item(84, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(85, 1).Value = '00:04:45.950'
$ws.Cells.Item(85, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(85, 3).Value = 'Ball Receipt'
$ws.Cells.Item(85, 4).Value = 4
$ws.Cells.Item(85, 5).Value = 33
$ws.Cells.Item(85, 6).Value = 74
$ws.Cells.Item(85, 12).Value = 'Adrian Pombo'
$ws.Cells.Item(85, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(86, 1).Value = '00:04:46.850'
$ws.Cells.Item(86, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(86, 3).Value = 'Pass'
$ws.Cells.Item(86, 4).Value = 4
$ws.Cells.Item(86, 5).Value = 30
$ws.Cells.Item(86, 6).Value = 74
$ws.Cells.Item(86, 7).Value = 15
$ws.Cells.Item(86, 8).Value = 50
$ws.Cells.Item(86, 9).Value = 'Alejandro Charro'
$ws.Cells.Item(86, 10).Value = 'Complete'
$ws.Cells.Item(86, 11).Value = 'High Pass'
$ws.Cells.Item(86, 12).Value = 'Adrian Pombo'
$ws.Cells.Item(86, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(87, 1).Value = '00:04:48.350'
$ws.Cells.Item(87, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(87, 3).Value = 'Miscontrol'
$ws.Cells.Item(87, 4).Value = 4
$ws.Cells.Item(87, 5).Value = 15
$ws.Cells.Item(87, 6).Value = 50
$ws.Cells.Item(87, 12).Value = 'Alejandro Charro'
$ws.Cells.Item(87, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(88, 1).Value = '00:04:49.850'
$ws.Cells.Item(88, 2).Value = 'C.D. Retiro Sur'
$ws.Cells.Item(88, 3).Value = 'Goal Keeper'
$ws.Cells.Item(88, 4).Value = 4
$ws.Cells.Item(88, 5).Value = 4
$ws.Cells.Item(88, 6).Value = 47
$ws.Cells.Item(88, 12).Value = 'Ignacio Salas'
$ws.Cells.Item(88, 13).Value = 'C.D. Retiro Sur'

$ws.Cells.Item(89, 1).Value = '00:05:09.250'
$ws.Cells.Item(89, 2).Value = 'C.D. Retiro Sur'
$ws.Cells.Item(89, 3).Value = 'Pass'
$ws.Cells.Item(89, 4).Value = 5
$ws.Cells.Item(89, 5).Value = 21
$ws.Cells.Item(89, 6).Value = 42
$ws.Cells.Item(89, 7).Value = 81
$ws.Cells.Item(89, 8).Value = 60
$ws.Cells.Item(89, 10).Value = 'Incomplete'
$ws.Cells.Item(89, 11).Value = 'High Pass'
$ws.Cells.Item(89, 12).Value = 'Ignacio Salas'
$ws.Cells.Item(89, 13).Value = 'C.D. Retiro Sur'

$ws.Cells.Item(90, 1).Value = '00:05:12.250'
$ws.Cells.Item(90, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(90, 3).Value = 'Ball Receipt'
$ws.Cells.Item(90, 4).Value = 5
$ws.Cells.Item(90, 5).Value = 81
$ws.Cells.Item(90, 6).Value = 60
$ws.Cells.Item(90, 12).Value = 'Adrian Pombo'
$ws.Cells.Item(90, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(91, 1).Value = '00:05:14.250'
$ws.Cells.Item(91, 2).Value = 'C.D. Retiro Sur'
$ws.Cells.Item(91, 3).Value = 'Ball Recovery'
$ws.Cells.Item(91, 4).Value = 5
$ws.Cells.Item(91, 5).Value = 83
$ws.Cells.Item(91, 6).Value = 58
$ws.Cells.Item(91, 12).Value = 'Juan Camilo'
$ws.Cells.Item(91, 13).Value = 'C.D. Retiro Sur'

$ws.Cells.Item(92, 1).Value = '00:05:15.290'
$ws.Cells.Item(92, 2).Value = 'C.D. Retiro Sur'
$ws.Cells.Item(92, 3).Value = 'Pass'
$ws.Cells.Item(92, 4).Value = 5
$ws.Cells.Item(92, 5).Value = 83
$ws.Cells.Item(92, 6).Value = 58
$ws.Cells.Item(92, 7).Value = 102
$ws.Cells.Item(92, 8).Value = 62
$ws.Cells.Item(92, 9).Value = 'Raúl Angullo'
$ws.Cells.Item(92, 10).Value = 'Complete'
$ws.Cells.Item(92, 11).Value = 'High Pass'
$ws.Cells.Item(92, 12).Value = 'Juan Camilo'
$ws.Cells.Item(92, 13).Value = 'C.D. Retiro Sur'

$ws.Cells.Item(93, 1).Value = '00:05:17.290'
$ws.Cells.Item(93, 2).Value = 'C.D. Retiro Sur'
$ws.Cells.Item(93, 3).Value = 'Ball Receipt'
$ws.Cells.Item(93, 4).Value = 5
$ws.Cells.Item(93, 5).Value = 106
$ws.Cells.Item(93, 6).Value = 63
$ws.Cells.Item(93, 12).Value = 'Raúl Angullo'
$ws.Cells.Item(93, 13).Value = 'C.D. Retiro Sur'

$ws.Cells.Item(94, 1).Value = '00:05:19.800'
$ws.Cells.Item(94, 2).Value = 'C.D. Retiro Sur'
$ws.Cells.Item(94, 3).Value = 'Misdribble'
$ws.Cells.Item(94, 4).Value = 5
$ws.Cells.Item(94, 5).Value = 108
$ws.Cells.Item(94, 6).Value = 65
$ws.Cells.Item(94, 12).Value = 'Raúl Angullo'
$ws.Cells.Item(94, 13).Value = 'C.D. Retiro Sur'

$ws.Cells.Item(95, 1).Value = '00:05:20.800'
$ws.Cells.Item(95, 2).Value = 'C.D. Retiro Sur'
$ws.Cells.Item(95, 3).Value = 'Foul Committed'
$ws.Cells.Item(95, 4).Value = 5
$ws.Cells.Item(95, 5).Value = 108
$ws.Cells.Item(95, 6).Value = 65
$ws.Cells.Item(95, 12).Value = 'Raúl Angullo'
$ws.Cells.Item(95, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(96, 1).Value = '00:05:33.700'
$ws.Cells.Item(96, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(96, 3).Value = 'Pass'
$ws.Cells.Item(96, 4).Value = 5
$ws.Cells.Item(96, 5).Value = 108
$ws.Cells.Item(96, 6).Value = 65
$ws.Cells.Item(96, 7).Value = 114
$ws.Cells.Item(96, 8).Value = 49
$ws.Cells.Item(96, 9).Value = 'Oscar Ponce'
$ws.Cells.Item(96, 10).Value = 'Complete'
$ws.Cells.Item(96, 11).Value = 'Ground Pass'
$ws.Cells.Item(96, 12).Value = 'Ángel Jesús'
$ws.Cells.Item(96, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(97, 1).Value = '00:05:35.940'
$ws.Cells.Item(97, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(97, 3).Value = 'Ball Receipt'
$ws.Cells.Item(97, 4).Value = 5
$ws.Cells.Item(97, 5).Value = 114
$ws.Cells.Item(97, 6).Value = 49
$ws.Cells.Item(97, 12).Value = 'Oscar Ponce'
$ws.Cells.Item(97, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(98, 1).Value = '00:05:37.110'
$ws.Cells.Item(98, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(98, 3).Value = 'Pass'
$ws.Cells.Item(98, 4).Value = 5
$ws.Cells.Item(98, 5).Value = 114
$ws.Cells.Item(98, 6).Value = 46
$ws.Cells.Item(98, 7).Value = 112
$ws.Cells.Item(98, 8).Value = 22
$ws.Cells.Item(98, 9).Value = 'Pablo Escribano'
$ws.Cells.Item(98, 10).Value = 'Complete'
$ws.Cells.Item(98, 11).Value = 'Ground Pass'
$ws.Cells.Item(98, 12).Value = 'Oscar Ponce'
$ws.Cells.Item(98, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(99, 1).Value = '00:05:38.680'
$ws.Cells.Item(99, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(99, 3).Value = 'Ball Receipt'
$ws.Cells.Item(99, 4).Value = 5
$ws.Cells.Item(99, 5).Value = 112
$ws.Cells.Item(99, 6).Value = 22
$ws.Cells.Item(99, 12).Value = 'Pablo Escribano'
$ws.Cells.Item(99, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(100, 1).Value = '00:05:40.680'
$ws.Cells.Item(100, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(100, 3).Value = 'Pass'
$ws.Cells.Item(100, 4).Value = 5
$ws.Cells.Item(100, 5).Value = 112
$ws.Cells.Item(100, 6).Value = 22
$ws.Cells.Item(100, 7).Value = 95
$ws.Cells.Item(100, 8).Value = 2
$ws.Cells.Item(100, 9).Value = 'Carlos Enrique'
$ws.Cells.Item(100, 10).Value = 'Complete'
$ws.Cells.Item(100, 11).Value = 'Ground Pass'
$ws.Cells.Item(100, 12).Value = 'Pablo Escribano'
$ws.Cells.Item(100, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(101, 1).Value = '00:05:42.680'
$ws.Cells.Item(101, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(101, 3).Value = 'Ball Receipt'
$ws.Cells.Item(101, 4).Value = 5
$ws.Cells.Item(101, 5).Value = 95
$ws.Cells.Item(101, 6).Value = 2
$ws.Cells.Item(101, 12).Value = 'Carlos Enrique'
$ws.Cells.Item(101, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(102, 1).Value = '00:05:43.680'
$ws.Cells.Item(102, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(102, 3).Value = 'Pass'
$ws.Cells.Item(102, 4).Value = 5
$ws.Cells.Item(102, 5).Value = 95
$ws.Cells.Item(102, 6).Value = 2
$ws.Cells.Item(102, 7).Value = 102
$ws.Cells.Item(102, 8).Value = 24
$ws.Cells.Item(102, 9).Value = 'Pablo Escribano'
$ws.Cells.Item(102, 10).Value = 'Complete'
$ws.Cells.Item(102, 11).Value = 'Ground Pass'
$ws.Cells.Item(102, 12).Value = 'Carlos Enrique'
$ws.Cells.Item(102, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(103, 1).Value = '00:05:45.400'
$ws.Cells.Item(103, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(103, 3).Value = 'Ball Receipt'
$ws.Cells.Item(103, 4).Value = 5
$ws.Cells.Item(103, 5).Value = 102
$ws.Cells.Item(103, 6).Value = 24
$ws.Cells.Item(103, 12).Value = 'Pablo Escribano'
$ws.Cells.Item(103, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(104, 1).Value = '00:05:46.630'
$ws.Cells.Item(104, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(104, 3).Value = 'Pass'
$ws.Cells.Item(104, 4).Value = 5
$ws.Cells.Item(104, 5).Value = 102
$ws.Cells.Item(104, 6).Value = 24
$ws.Cells.Item(104, 7).Value = 100
$ws.Cells.Item(104, 8).Value = 64
$ws.Cells.Item(104, 9).Value = 'Ángel Jesús'
$ws.Cells.Item(104, 10).Value = 'Complete'
$ws.Cells.Item(104, 11).Value = 'Ground Pass'
$ws.Cells.Item(104, 12).Value = 'Pablo Escribano'
$ws.Cells.Item(104, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(105, 1).Value = '00:05:49.640'
$ws.Cells.Item(105, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(105, 3).Value = 'Ball Receipt'
$ws.Cells.Item(105, 4).Value = 5
$ws.Cells.Item(105, 5).Value = 100
$ws.Cells.Item(105, 6).Value = 64
$ws.Cells.Item(105, 12).Value = 'Ángel Jesús'
$ws.Cells.Item(105, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(106, 1).Value = '00:05:51.640'
$ws.Cells.Item(106, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(106, 3).Value = 'Pass'
$ws.Cells.Item(106, 4).Value = 5
$ws.Cells.Item(106, 5).Value = 90
$ws.Cells.Item(106, 6).Value = 65
$ws.Cells.Item(106, 7).Value = 80
$ws.Cells.Item(106, 8).Value = 79
$ws.Cells.Item(106, 9).Value = 'Adrian Pombo'
$ws.Cells.Item(106, 10).Value = 'Complete'
$ws.Cells.Item(106, 11).Value = 'Ground Pass'
$ws.Cells.Item(106, 12).Value = 'Ángel Jesús'
$ws.Cells.Item(106, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(107, 1).Value = '00:05:53.640'
$ws.Cells.Item(107, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(107, 3).Value = 'Ball Receipt'
$ws.Cells.Item(107, 4).Value = 5
$ws.Cells.Item(107, 5).Value = 80
$ws.Cells.Item(107, 6).Value = 79
$ws.Cells.Item(107, 12).Value = 'Adrian Pombo'
$ws.Cells.Item(107, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(108, 1).Value = '00:05:55.310'
$ws.Cells.Item(108, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(108, 3).Value = 'Pass'
$ws.Cells.Item(108, 4).Value = 5
$ws.Cells.Item(108, 5).Value = 80
$ws.Cells.Item(108, 6).Value = 76
$ws.Cells.Item(108, 7).Value = 86
$ws.Cells.Item(108, 8).Value = 77
$ws.Cells.Item(108, 9).Value = 'Ángel Jesús'
$ws.Cells.Item(108, 10).Value = 'Complete'
$ws.Cells.Item(108, 11).Value = 'Ground Pass'
$ws.Cells.Item(108, 12).Value = 'Adrian Pombo'
$ws.Cells.Item(108, 13).Value = 'Escuela Dep. Moratalaz ''D'''

$ws.Cells.Item(109, 1).Value = '00:05:56.510'
$ws.Cells.Item(109, 2).Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Cells.Item(109, 3).Value = 'Ball Receipt'
$ws.Cells.Item(109, 4).Value = 5
$ws.Cells.Item(109, 5).Value = 86
$ws.Cells.Item(109, 6).Value = 77
$ws.Cells.Item(109, 12).Value = 'Ángel Jesús'
$ws.Cells.Item(109, 13).Value = 'Escuela Dep. Moratalaz ''D'''

# Restore the view to where the author left it: scrolled near the bottom of
# the newly appended block, with Q101 as the active cell.
$ws.Range("Q101").Select()
